$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.319.35"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.798.79"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'595.56"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'168.14"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "3.797.75"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").Value = "'36.26"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "4.438.80"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "3.837.83"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "68.390.76"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'17.85"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").Value = "'7.00"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D21").Value = "'10.67"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "'463.89"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "'0.699"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  +6.55%  "
$ws.Range("D25").Value = "'83.90"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Value = "'11.90"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "'30.02"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "'7.25"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "'2.16"
$ws.Range("E33").Value = "  -4.30%  "
$ws.Range("D34").Value = "'9.13"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "3.751.37"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "'3.51"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").Value = "'0.994"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.301"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").Value = "'43.64"
$ws.Range("E45").Value = "  +8.41%  "
$ws.Range("D46").Value = "'46.95"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").Value = "'8.40"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'146.96"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "'389.76"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000265"
$ws.Range("E51").Value = "  +4.03%  "

Write-Output "Applied 89 cell updates"
